$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-26 down to 13-27.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with data (matching the constant columns used by all
# other data rows, plus the new record's specific values).
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = [DateTime]::FromOADate(44671)
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112001
$ws.Cells.Item(12, 7).Value = "Berenjena"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 160
$ws.Cells.Item(12, 11).Value = 6000
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 13).Value = 6500
$ws.Cells.Item(12, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 108
$ws.Cells.Item(12, 17).Value = 60
$ws.Cells.Item(12, 18).Value = "Hortaliza"
